$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("AI4").Value = 0.413
$ws.Range("AJ4").Value = 0.094
$ws.Range("AK4").Value = 0.306
$ws.Range("AU4").Value = 0.246
$ws.Range("AW4").Value = 0.16
$ws.Range("BA4").Value = 2.03
$ws.Range("BB4").Value = 0.146
$ws.Range("BC4").Value = 0.382
$ws.Range("BG4").Value = 0.715
$ws.Range("BH4").Value = 0.144
$ws.Range("BI4").Value = 0.379
$ws.Range("BM4").Value = 0.744
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.677
$ws.Range("BQ4").Value = 0.759
$ws.Range("E4").Value = 0.476
$ws.Range("F4").Value = 0.053
$ws.Range("G4").Value = 0.231
$ws.Range("N4").Value = 0.489
$ws.Range("O4").Value = 0.058
$ws.Range("P4").Value = 0.241
$ws.Range("Q4").Value = 0.054
$ws.Range("R4").Value = 0.036
$ws.Range("S4").Value = 0.19
$ws.Range("W4").Value = 0.376
$ws.Range("X4").Value = 0.107
$ws.Range("Y4").Value = 0.328
$ws.Range("AI5").Value = 0.413
$ws.Range("AJ5").Value = 0.09
$ws.Range("AK5").Value = 0.3
$ws.Range("AU5").Value = 0.459
$ws.Range("AV5").Value = 0.075
$ws.Range("AW5").Value = 0.275
$ws.Range("BA5").Value = 1.306
$ws.Range("BG5").Value = 0.376
$ws.Range("BH5").Value = 0.05
$ws.Range("BI5").Value = 0.224
$ws.Range("BM5").Value = 0.53
$ws.Range("BN5").Value = 0.047
$ws.Range("BO5").Value = 0.217
$ws.Range("BP5").Value = 0.435
$ws.Range("BQ5").Value = 0.458
$ws.Range("E5").Value = 0.604
$ws.Range("F5").Value = 0.06
$ws.Range("G5").Value = 0.244
$ws.Range("N5").Value = 0.734
$ws.Range("O5").Value = 0.066
$ws.Range("P5").Value = 0.258
$ws.Range("Q5").Value = 0.035
$ws.Range("R5").Value = 0.015
$ws.Range("S5").Value = 0.122
$ws.Range("W5").Value = 0.343
$ws.Range("X5").Value = 0.099
$ws.Range("Y5").Value = 0.315
$ws.Range("AI6").Value = 0.413
$ws.Range("AU6").Value = 0.32
$ws.Range("BA6").Value = 1.582
$ws.Range("BG6").Value = 0.493
$ws.Range("BM6").Value = 0.619
$ws.Range("BP6").Value = 0.527
$ws.Range("BQ6").Value = 0.569
$ws.Range("E6").Value = 0.532
$ws.Range("N6").Value = 0.587
$ws.Range("Q6").Value = 0.042
$ws.Range("W6").Value = 0.359
$ws.Range("AI7").Value = 0.413
$ws.Range("AU7").Value = 0.391
$ws.Range("BA7").Value = 1.402
$ws.Range("BG7").Value = 0.415
$ws.Range("BM7").Value = 0.5620000000000001
$ws.Range("BP7").Value = 0.467
$ws.Range("BQ7").Value = 0.497
$ws.Range("E7").Value = 0.573
$ws.Range("N7").Value = 0.667
$ws.Range("Q7").Value = 0.038
$ws.Range("W7").Value = 0.349
$ws.Range("AI8").Value = 0.483
$ws.Range("AJ8").Value = 0.138
$ws.Range("AK8").Value = 0.371
$ws.Range("AU8").Value = 0.401
$ws.Range("AV8").Value = 0.083
$ws.Range("AW8").Value = 0.289
$ws.Range("BA8").Value = 1.768
$ws.Range("BB8").Value = 0.111
$ws.Range("BC8").Value = 0.333
$ws.Range("BG8").Value = 0.5639999999999999
$ws.Range("BH8").Value = 0.111
$ws.Range("BI8").Value = 0.333
$ws.Range("BM8").Value = 0.6820000000000001
$ws.Range("BN8").Value = 0.06
$ws.Range("BO8").Value = 0.246
$ws.Range("BP8").Value = 0.589
$ws.Range("BQ8").Value = 0.626
$ws.Range("E8").Value = 0.703
$ws.Range("F8").Value = 0.074
$ws.Range("G8").Value = 0.271
$ws.Range("N8").Value = 0.822
$ws.Range("O8").Value = 0.045
$ws.Range("P8").Value = 0.212
$ws.Range("Q8").Value = 0.039
$ws.Range("W8").Value = 0.417
$ws.Range("X8").Value = 0.121
$ws.Range("Y8").Value = 0.347
$ws.Range("AI9").Value = 0.439
$ws.Range("AJ9").Value = 0.246
$ws.Range("AK9").Value = 0.496
$ws.Range("BA9").Value = 1.732
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.61
$ws.Range("BH9").Value = 0.238
$ws.Range("BI9").Value = 0.488
$ws.Range("BM9").Value = 0.659
$ws.Range("BN9").Value = 0.225
$ws.Range("BO9").Value = 0.474
$ws.Range("BP9").Value = 0.577
$ws.Range("BQ9").Value = 0.618
$ws.Range("E9").Value = 0.659
$ws.Range("F9").Value = 0.225
$ws.Range("G9").Value = 0.474
$ws.Range("N9").Value = 0.756
$ws.Range("O9").Value = 0.184
$ws.Range("P9").Value = 0.429
$ws.Range("W9").Value = 0.317
$ws.Range("X9").Value = 0.217
$ws.Range("Y9").Value = 0.465
$ws.Range("AI10").Value = 0.512
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.39
$ws.Range("AV10").Value = 0.238
$ws.Range("AW10").Value = 0.488
$ws.Range("BA10").Value = 2.195
$ws.Range("BB10").Value = 0.217
$ws.Range("BC10").Value = 0.465
$ws.Range("BG10").Value = 0.6830000000000001
$ws.Range("BH10").Value = 0.217
$ws.Range("BI10").Value = 0.465
$ws.Range("BM10").Value = 0.829
$ws.Range("BN10").Value = 0.142
$ws.Range("BO10").Value = 0.376
$ws.Range("BP10").Value = 0.732
$ws.Range("BQ10").Value = 0.765
$ws.Range("E10").Value = 0.805
$ws.Range("F10").Value = 0.157
$ws.Range("G10").Value = 0.396
$ws.Range("N10").Value = 0.951
$ws.Range("O10").Value = 0.046
$ws.Range("P10").Value = 0.215
$ws.Range("W10").Value = 0.537
$ws.Range("AI11").Value = 0.585
$ws.Range("AJ11").Value = 0.243
$ws.Range("AK11").Value = 0.493
$ws.Range("AU11").Value = 0.5610000000000001
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.195
$ws.Range("BB11").Value = 0.217
$ws.Range("BC11").Value = 0.465
$ws.Range("BG11").Value = 0.6830000000000001
$ws.Range("BH11").Value = 0.217
$ws.Range("BI11").Value = 0.465
$ws.Range("BM11").Value = 0.829
$ws.Range("BN11").Value = 0.142
$ws.Range("BO11").Value = 0.376
$ws.Range("BP11").Value = 0.732
$ws.Range("BQ11").Value = 0.772
$ws.Range("E11").Value = 0.854
$ws.Range("F11").Value = 0.125
$ws.Range("G11").Value = 0.353
$ws.Range("N11").Value = 0.951
$ws.Range("O11").Value = 0.046
$ws.Range("P11").Value = 0.215
$ws.Range("W11").Value = 0.537
$ws.Range("AI12").Value = 1.583
$ws.Range("AJ12").Value = 1.493
$ws.Range("AK12").Value = 1.222
$ws.Range("AU12").Value = 2.88
$ws.Range("AV12").Value = 3.466
$ws.Range("AW12").Value = 1.862
$ws.Range("BA12").Value = 3.807
$ws.Range("BB12").Value = 0.459
$ws.Range("BC12").Value = 0.678
$ws.Range("BG12").Value = 1.143
$ws.Range("BH12").Value = 0.194
$ws.Range("BI12").Value = 0.44
$ws.Range("BM12").Value = 1.235
$ws.Range("BN12").Value = 0.239
$ws.Range("BO12").Value = 0.489
$ws.Range("BP12").Value = 1.269
$ws.Range("BQ12").Value = 1.253
$ws.Range("E12").Value = 1.429
$ws.Range("F12").Value = 0.873
$ws.Range("G12").Value = 0.9350000000000001
$ws.Range("N12").Value = 1.256
$ws.Range("O12").Value = 0.293
$ws.Range("P12").Value = 0.542
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.432
$ws.Range("Y12").Value = 0.657
$ws.Range("AI13").Value = 1.159
$ws.Range("AJ13").Value = 0.312
$ws.Range("AK13").Value = 0.5580000000000001
$ws.Range("AU13").Value = 2.048
$ws.Range("AV13").Value = 0.344
$ws.Range("AW13").Value = 0.587
$ws.Range("BA13").Value = 2.2
$ws.Range("BB13").Value = 0.284
$ws.Range("BC13").Value = 0.533
$ws.Range("BG13").Value = 0.544
$ws.Range("BH13").Value = 0.051
$ws.Range("BI13").Value = 0.226
$ws.Range("BM13").Value = 0.798
$ws.Range("BN13").Value = 0.161
$ws.Range("BO13").Value = 0.402
$ws.Range("BP13").Value = 0.733
$ws.Range("BQ13").Value = 0.669
$ws.Range("E13").Value = 1.425
$ws.Range("F13").Value = 0.298
$ws.Range("G13").Value = 0.546
$ws.Range("N13").Value = 1.73
$ws.Range("O13").Value = 0.476
$ws.Range("P13").Value = 0.6899999999999999
$ws.Range("W13").Value = 0.985
$ws.Range("X13").Value = 0.199
$ws.Range("Y13").Value = 0.446
